$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Make 윤다은 the active/selected sheet (was 문준범 before)
$ws.Activate()

# Row 5: increase wrap height
$ws.Rows.Item(5).RowHeight = 153

# Row 6: fill in content/result/date columns
$ws.Range("B6").Value = "개인정보 수정 구현 (디비에서 정보를 가져와서 보여준후 수정된 내용 다시 반영)"
$ws.Range("E6").Value = "완료"

# Row 5 "문제점" text update (replaces old shared string so it is garbage-collected
# and "환경구축완료" shifts down to index 35, matching E3)
$ws.Range("F5").Value = "신고할 때 입력된 내용에대한 처리가 없음-> 5회이상 신고당한 유저 차단, 신고된 내용에 대해 경고횟수 증가시키는 부분 구현필요-> 구현"

$ws.Range("F6").Value = "현재 비밀번호 확인하는 부분이 없음"
$ws.Range("D6").Value = 43618

# Row 7: add assigned date
$ws.Range("C7").Value = 43611

# Row 8
$ws.Range("A8").Value = "게시글 등록기능구현"
$ws.Range("C8").Value = 43618

# Row 9
$ws.Range("A9").Value = "게시글 보기, 검색"
$ws.Range("C9").Value = 43618

# Row 10
$ws.Range("A10").Value = "신청내역보기 및 승인"
$ws.Range("C10").Value = 43618

# Row 11
$ws.Rows.Item(11).RowHeight = 34.5
$ws.Range("A11").Value = "이전까지 신청한 히스토리보기"
$ws.Range("C11").Value = 43618

# Row 12
$ws.Range("A12").Value = "사용자 개인 페이지구현"
$ws.Range("C12").Value = 43618

# Update the selection shown on the now-active 윤다은 sheet
$ws.Range("B10").Select()
